$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "The time is almost the same..." note in N1
$ws.Range("N1").ClearContents()

# --- Left table (A:C) - Number Of Processes / Time / Number Of Prime Numbers ---
# New rows for 1 and 2 processes are inserted before the existing data, and
# the old "15000" row is dropped - 8 data rows total (rows 2-9).
$leftData = @(
  @(1,    0.03125,    84000000),
  @(2,    0.046875,   84000000),
  @(10,   0.125,      84000000),
  @(100,  1.234375,   84000000),
  @(1000, 12.859375,  84000000),
  @(2000, 25,         84000000),
  @(5000, 63.828125,  84000000),
  @(10000,134.0625,   84000000)
)
for ($i = 0; $i -lt $leftData.Length; $i++) {
  $r = $i + 2
  $row = $leftData[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- Right table (H:J) - Number Of Processes / Time / Number Of Prime Numbers ---
# New rows for 1 and 2 processes inserted before existing data, plus a new
# trailing row for 20000 processes - 9 data rows total (rows 2-10).
$rightData = @(
  @(1,    0,          1000000),
  @(2,    0,          1000000),
  @(10,   0,          1000000),
  @(100,  0.109375,   1000000),
  @(1000, 0.8125,     1000000),
  @(2000, 1.859375,   1000000),
  @(5000, 6.640625,   1000000),
  @(10000,18.109375,  1000000),
  @(20000,60.375,     1000000)
)
for ($i = 0; $i -lt $rightData.Length; $i++) {
  $r = $i + 2
  $row = $rightData[$i]
  $ws.Cells.Item($r, 8).Value = $row[0]
  $ws.Cells.Item($r, 9).Value = $row[1]
  $ws.Cells.Item($r, 10).Value = $row[2]
}

# Update the saved selection (matches the author's last-saved cursor position)
$ws.Range("G20").Select()

# --- Chart 1 (left table A:B) ---
$chart1 = $ws.ChartObjects().Item(1).Chart
$series1 = $chart1.SeriesCollection(1)
$series1.XValues = "=Arkusz1!`$A`$3:`$A`$10"
$series1.Values  = "=Arkusz1!`$B`$3:`$B`$10"

$tl1 = $series1.Trendlines().Item(1)
$tl1.Type = -4132
$tl1.Order = $null

$yAxis1 = $chart1.Axes(2)
$yAxis1.MaximumScale = 150
$yAxis1.MinimumScale = 0

# --- Chart 2 (right table H:I) ---
$chart2 = $ws.ChartObjects().Item(2).Chart
$series2 = $chart2.SeriesCollection(1)
$series2.XValues = "=Arkusz1!`$H`$2:`$H`$10"
$series2.Values  = "=Arkusz1!`$I`$2:`$I`$10"
